$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1463976666666667
$ws.Range("H2").Value = 0.439193
$ws.Range("I2").Value = 0.1157910139257259
$ws.Range("J2").Value = 0.115791013925726
$ws.Range("O2").Value = 0.06522509891308133
$ws.Range("P2").Value = 0.06522509891308133
$ws.Range("Q2").Value = 0.02998448689755556
$ws.Range("R2").Value = 0.269860382078
$ws.Range("S2").Value = 0.007552480336551452
$ws.Range("T2").Value = 0.007552480336551453

$ws.Range("G3").Value = 0.1463976666666667
$ws.Range("H3").Value = 0.439193
$ws.Range("I3").Value = 0.1157910139257259
$ws.Range("J3").Value = 0.115791013925726
$ws.Range("M3").Value = 0.6481333333333333
$ws.Range("N3").Value = 1.9444
$ws.Range("O3").Value = 0.2064033004146749
$ws.Range("P3").Value = 0.2064033004146749
$ws.Range("Q3").Value = 0.09488520768888889
$ws.Range("R3").Value = 0.8539668692
$ws.Range("S3").Value = 0.02389964743263142
$ws.Range("T3").Value = 0.02389964743263142

$ws.Range("G4").Value = 0.1463976666666667
$ws.Range("H4").Value = 0.439193
$ws.Range("I4").Value = 0.1157910139257259
$ws.Range("J4").Value = 0.115791013925726
$ws.Range("M4").Value = 1.888205
$ws.Range("N4").Value = 5.664615
$ws.Range("O4").Value = 0.6013141491351952
$ws.Range("P4").Value = 0.6013141491351952
$ws.Range("Q4").Value = 0.2764288061883333
$ws.Range("R4").Value = 2.487859255695
$ws.Range("S4").Value = 0.06962677501624943
$ws.Range("T4").Value = 0.06962677501624943

$ws.Range("G5").Value = 0.1463976666666667
$ws.Range("H5").Value = 0.439193
$ws.Range("I5").Value = 0.1157910139257259
$ws.Range("J5").Value = 0.115791013925726
$ws.Range("M5").Value = 0.398977
$ws.Range("N5").Value = 1.196931
$ws.Range("O5").Value = 0.1270574515370486
$ws.Range("P5").Value = 0.1270574515370486
$ws.Range("Q5").Value = 0.05840930185366667
$ws.Range("R5").Value = 0.525683716683
$ws.Range("S5").Value = 0.01471211114029364
$ws.Range("T5").Value = 0.01471211114029364

$ws.Range("I6").Value = 0.4041732358198567
$ws.Range("J6").Value = 0.4041732358198568
$ws.Range("O6").Value = 0.06522509891308133
$ws.Range("P6").Value = 0.06522509891308133
$ws.Range("S6").Value = 0.0263622392843703
$ws.Range("T6").Value = 0.0263622392843703

$ws.Range("I7").Value = 0.4041732358198567
$ws.Range("J7").Value = 0.4041732358198568
$ws.Range("M7").Value = 0.6481333333333333
$ws.Range("N7").Value = 1.9444
$ws.Range("O7").Value = 0.2064033004146749
$ws.Range("P7").Value = 0.2064033004146749
$ws.Range("Q7").Value = 0.3312006702666667
$ws.Range("R7").Value = 2.9808060324
$ws.Range("S7").Value = 0.08342268981249713
$ws.Range("T7").Value = 0.08342268981249713

$ws.Range("I8").Value = 0.4041732358198567
$ws.Range("J8").Value = 0.4041732358198568
$ws.Range("M8").Value = 1.888205
$ws.Range("N8").Value = 5.664615
$ws.Range("O8").Value = 0.6013141491351952
$ws.Range("P8").Value = 0.6013141491351952
$ws.Range("Q8").Value = 0.9648859724350002
$ws.Range("R8").Value = 8.683973751915
$ws.Range("S8").Value = 0.2430350854002357
$ws.Range("T8").Value = 0.2430350854002358

$ws.Range("I9").Value = 0.4041732358198567
$ws.Range("J9").Value = 0.4041732358198568
$ws.Range("M9").Value = 0.398977
$ws.Range("N9").Value = 1.196931
$ws.Range("O9").Value = 0.1270574515370486
$ws.Range("P9").Value = 0.1270574515370486
$ws.Range("Q9").Value = 0.203880039839
$ws.Range("R9").Value = 1.834920358551
$ws.Range("S9").Value = 0.05135322132275354
$ws.Range("T9").Value = 0.05135322132275355

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1328766666666667
$ws.Range("H10").Value = 0.39863
$ws.Range("I10").Value = 0.1050967840589721
$ws.Range("J10").Value = 0.1050967840589721
$ws.Range("O10").Value = 0.06522509891308133
$ws.Range("P10").Value = 0.06522509891308133
$ws.Range("Q10").Value = 0.02721517877555556
$ws.Range("R10").Value = 0.24493660898
$ws.Range("S10").Value = 0.006854948135693204
$ws.Range("T10").Value = 0.006854948135693204

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1328766666666667
$ws.Range("H11").Value = 0.39863
$ws.Range("I11").Value = 0.1050967840589721
$ws.Range("J11").Value = 0.1050967840589721
$ws.Range("M11").Value = 0.6481333333333333
$ws.Range("N11").Value = 1.9444
$ws.Range("O11").Value = 0.2064033004146749
$ws.Range("P11").Value = 0.2064033004146749
$ws.Range("Q11").Value = 0.08612179688888889
$ws.Range("R11").Value = 0.7750961719999999
$ws.Range("S11").Value = 0.02169232309274024
$ws.Range("T11").Value = 0.02169232309274023

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1328766666666667
$ws.Range("H12").Value = 0.39863
$ws.Range("I12").Value = 0.1050967840589721
$ws.Range("J12").Value = 0.1050967840589721
$ws.Range("M12").Value = 1.888205
$ws.Range("N12").Value = 5.664615
$ws.Range("O12").Value = 0.6013141491351952
$ws.Range("P12").Value = 0.6013141491351952
$ws.Range("Q12").Value = 0.2508983863833333
$ws.Range("R12").Value = 2.25808547745
$ws.Range("S12").Value = 0.06319618328326615
$ws.Range("T12").Value = 0.06319618328326615

$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1328766666666667
$ws.Range("H13").Value = 0.39863
$ws.Range("I13").Value = 0.1050967840589721
$ws.Range("J13").Value = 0.1050967840589721
$ws.Range("M13").Value = 0.398977
$ws.Range("N13").Value = 1.196931
$ws.Range("O13").Value = 0.1270574515370486
$ws.Range("P13").Value = 0.1270574515370486
$ws.Range("Q13").Value = 0.05301473383666666
$ws.Range("R13").Value = 0.4771326045299999
$ws.Range("S13").Value = 0.0133533295472725
$ws.Range("T13").Value = 0.0133533295472725

$ws.Range("G14").Value = 0.180116
$ws.Range("H14").Value = 0.5403480000000001
$ws.Range("I14").Value = 0.1424600182442301
$ws.Range("J14").Value = 0.1424600182442301
$ws.Range("O14").Value = 0.06522509891308133
$ws.Range("P14").Value = 0.06522509891308133
$ws.Range("Q14").Value = 0.03689051857866667
$ws.Range("R14").Value = 0.332014667208
$ws.Range("S14").Value = 0.009291968781139281
$ws.Range("T14").Value = 0.009291968781139281

$ws.Range("G15").Value = 0.180116
$ws.Range("H15").Value = 0.5403480000000001
$ws.Range("I15").Value = 0.1424600182442301
$ws.Range("J15").Value = 0.1424600182442301
$ws.Range("M15").Value = 0.6481333333333333
$ws.Range("N15").Value = 1.9444
$ws.Range("O15").Value = 0.2064033004146749
$ws.Range("P15").Value = 0.2064033004146749
$ws.Range("Q15").Value = 0.1167391834666667
$ws.Range("R15").Value = 1.0506526512
$ws.Range("S15").Value = 0.0294042179427439
$ws.Range("T15").Value = 0.0294042179427439

$ws.Range("G16").Value = 0.180116
$ws.Range("H16").Value = 0.5403480000000001
$ws.Range("I16").Value = 0.1424600182442301
$ws.Range("J16").Value = 0.1424600182442301
$ws.Range("M16").Value = 1.888205
$ws.Range("N16").Value = 5.664615
$ws.Range("O16").Value = 0.6013141491351952
$ws.Range("P16").Value = 0.6013141491351952
$ws.Range("Q16").Value = 0.34009593178
$ws.Range("R16").Value = 3.06086338602
$ws.Range("S16").Value = 0.08566322465631362
$ws.Range("T16").Value = 0.08566322465631362

$ws.Range("G17").Value = 0.180116
$ws.Range("H17").Value = 0.5403480000000001
$ws.Range("I17").Value = 0.1424600182442301
$ws.Range("J17").Value = 0.1424600182442301
$ws.Range("M17").Value = 0.398977
$ws.Range("N17").Value = 1.196931
$ws.Range("O17").Value = 0.1270574515370486
$ws.Range("P17").Value = 0.1270574515370486
$ws.Range("Q17").Value = 0.071862141332
$ws.Range("R17").Value = 0.6467592719880001
$ws.Range("S17").Value = 0.01810060686403332
$ws.Range("T17").Value = 0.01810060686403332

$ws.Range("G18").Value = 0.2939293333333333
$ws.Range("H18").Value = 0.881788
$ws.Range("I18").Value = 0.2324789479512151
$ws.Range("J18").Value = 0.2324789479512152
$ws.Range("O18").Value = 0.06522509891308133
$ws.Range("P18").Value = 0.06522509891308133
$ws.Range("Q18").Value = 0.06020123438311111
$ws.Range("R18").Value = 0.541811109448
$ws.Range("S18").Value = 0.01516346237532709
$ws.Range("T18").Value = 0.01516346237532709

$ws.Range("G19").Value = 0.2939293333333333
$ws.Range("H19").Value = 0.881788
$ws.Range("I19").Value = 0.2324789479512151
$ws.Range("J19").Value = 0.2324789479512152
$ws.Range("M19").Value = 0.6481333333333333
$ws.Range("N19").Value = 1.9444
$ws.Range("O19").Value = 0.2064033004146749
$ws.Range("P19").Value = 0.2064033004146749
$ws.Range("Q19").Value = 0.1905053985777778
$ws.Range("R19").Value = 1.7145485872
$ws.Range("S19").Value = 0.04798442213406223
$ws.Range("T19").Value = 0.04798442213406223

$ws.Range("G20").Value = 0.2939293333333333
$ws.Range("H20").Value = 0.881788
$ws.Range("I20").Value = 0.2324789479512151
$ws.Range("J20").Value = 0.2324789479512152
$ws.Range("M20").Value = 1.888205
$ws.Range("N20").Value = 5.664615
$ws.Range("O20").Value = 0.6013141491351952
$ws.Range("P20").Value = 0.6013141491351952
$ws.Range("Q20").Value = 0.5549988368466666
$ws.Range("R20").Value = 4.99498953162
$ws.Range("S20").Value = 0.1397928807791302
$ws.Range("T20").Value = 0.1397928807791303

$ws.Range("G21").Value = 0.2939293333333333
$ws.Range("H21").Value = 0.881788
$ws.Range("I21").Value = 0.2324789479512151
$ws.Range("J21").Value = 0.2324789479512152
$ws.Range("M21").Value = 0.398977
$ws.Range("N21").Value = 1.196931
$ws.Range("O21").Value = 0.1270574515370486
$ws.Range("P21").Value = 0.1270574515370486
$ws.Range("Q21").Value = 0.1172710436253333
$ws.Range("R21").Value = 1.055439392628
$ws.Range("S21").Value = 0.02953818266269555
$ws.Range("T21").Value = 0.02953818266269555
